$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column C, shifting the old
# column C (and the "Jun_10" header it carried) two places to the right
# so it lands in column E.
$ws.Range("C:D").EntireColumn.Insert()

# Match the column width/format that the (now shifted) neighbour column
# already used for every column in the newly widened block.
$ws.Range("C:E").EntireColumn.ColumnWidth = 7.14

# Refresh the date-column headers: two new (most-recent) dates up front,
# followed by the two pre-existing dates that are now one column further
# along the row.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"
$ws.Range("D1").Value = "Jun_13"

# Default rating ("UN") for every data row in the two new columns.
$ws.Range("C2:D27").Value = "UN"
